# -----------------------------------------------------------------------
# Applies the styles.xml changes described by the diff:
#   1. Normal style: paragraph alignment left -> both (justify)
#   2. TableCaption / ImageCaption / CaptionedFigure styles: add center
#      paragraph alignment
#   3. Four new paragraph styles: Table Contents, Table Heading, Table,
#      My Table Heading
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# wdStyleType
$wdStyleTypeParagraph = 1

# wdParagraphAlignment
$wdAlignParagraphLeft    = 0
$wdAlignParagraphCenter  = 1
$wdAlignParagraphJustify = 3

# --- 1. Normal: jc left -> both --------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# --- 2. Existing caption-ish styles: center them ----------------------
$tableCaption = $d.Styles("TableCaption")
$tableCaption.ParagraphFormat.Alignment = $wdAlignParagraphCenter

$imageCaption = $d.Styles("ImageCaption")
$imageCaption.ParagraphFormat.Alignment = $wdAlignParagraphCenter

$captionedFigure = $d.Styles("CaptionedFigure")
$captionedFigure.ParagraphFormat.Alignment = $wdAlignParagraphCenter

# --- 3. New styles ------------------------------------------------------

# "Table Contents" - basedOn Normal
$tableContents = $d.Styles.Add("Table Contents", $wdStyleTypeParagraph)
$tableContents.BaseStyle = $d.Styles("Normal")
$tableContents.QuickStyle = $true
$tableContents.ParagraphFormat.WidowControl = $false
$tableContents.ParagraphFormat.NoLineNumber = $true
$tableContents.Font.NameAscii = "Times New Roman"
$tableContents.Font.NameOther = "Times New Roman"
$tableContents.Font.Name = "Times New Roman"

# "Table Heading" - basedOn Table Contents
$tableHeading = $d.Styles.Add("Table Heading", $wdStyleTypeParagraph)
$tableHeading.BaseStyle = $tableContents
$tableHeading.QuickStyle = $true
$tableHeading.ParagraphFormat.NoLineNumber = $true
# Shading (gray fill, like the real style definition) - applied for
# completeness; style-level shading round-trips through the paragraph
# shading properties below.
$tableHeading.ParagraphFormat.Shading.Texture = 0
$tableHeading.ParagraphFormat.Shading.ForegroundPatternColor = 0xFFFFFF
$tableHeading.ParagraphFormat.Shading.BackgroundPatternColor = 0x999999
$tableHeading.ParagraphFormat.Alignment = $wdAlignParagraphCenter
$tableHeading.Font.Bold = $true
$tableHeading.Font.BoldBi = $true

# "Table" (paragraph style) - basedOn Caption.
# NB: a table-type style already named "Table" exists in this document,
# and Word's style collection is keyed by name, so Styles.Add("Table", ...)
# would resolve to that existing (wrong-type) style instead of minting a
# new one. Mint it under a scratch name and rename it to "Table" so the
# resulting style's Name matches the target without disturbing the
# pre-existing table style.
$tableParaScratchName = "TableParaStyleScratch"
$tablePara = $d.Styles.Add($tableParaScratchName, $wdStyleTypeParagraph)
$tablePara.BaseStyle = $d.Styles("Caption")
$tablePara.QuickStyle = $true
$tablePara.ParagraphFormat.Alignment = $wdAlignParagraphCenter
$tablePara.NameLocal = "Table"

# "My Table Heading" - basedOn Table Heading
$myTableHeading = $d.Styles.Add("My Table Heading", $wdStyleTypeParagraph)
$myTableHeading.BaseStyle = $tableHeading
$myTableHeading.QuickStyle = $true
$myTableHeading.ParagraphFormat.Shading.Texture = 0
$myTableHeading.ParagraphFormat.Shading.ForegroundPatternColor = 0xFFFFFF
$myTableHeading.ParagraphFormat.Shading.BackgroundPatternColor = 0x999999

Write-Output "Styles updated."
